# CIERRRE 27 NOV 2021
# Fill in the credit/remision data rows for NOVIEMBRE sheet (rows 9-26)
# and update the sheet view's selection to reflect where the user left off.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)

$ws.Cells.Item(9,1).Value = 44510
$ws.Cells.Item(9,4).Value = "CANCELADA"
$ws.Cells.Item(9,5).Value = 0

$ws.Cells.Item(10,1).Value = 44510
$ws.Cells.Item(10,4).Value = "COMERCIO CENTRAL "
$ws.Cells.Item(10,5).Value = 1333

$ws.Cells.Item(11,1).Value = 44511
$ws.Cells.Item(11,4).Value = "OBRADOR"
$ws.Cells.Item(11,5).Value = 150

$ws.Cells.Item(12,1).Value = 44511
$ws.Cells.Item(12,4).Value = "COMERCIO CENTRAL "
$ws.Cells.Item(12,5).Value = 2460

$ws.Cells.Item(13,1).Value = 44512
$ws.Cells.Item(13,4).Value = "OBRADOR"
$ws.Cells.Item(13,5).Value = 8923

$ws.Cells.Item(14,1).Value = 44512
$ws.Cells.Item(14,4).Value = "COMERCIO CENTRAL "
$ws.Cells.Item(14,5).Value = 47911

$ws.Cells.Item(15,1).Value = 44512
$ws.Cells.Item(15,4).Value = "CANCELADA"
$ws.Cells.Item(15,4).Font.Bold = $true
$ws.Cells.Item(15,4).Font.Size = 12
$ws.Cells.Item(15,4).Font.Color = 255
$ws.Cells.Item(15,5).Value = 0

$ws.Cells.Item(16,1).Value = 44512
$ws.Cells.Item(16,4).Value = "COMERCIO CENTRAL "
$ws.Cells.Item(16,5).Value = 622

$ws.Cells.Item(17,1).Value = 44512
$ws.Cells.Item(17,4).Value = "COMERCIO CENTRAL "
$ws.Cells.Item(17,5).Value = 10714

$ws.Cells.Item(18,1).Value = 44512
$ws.Cells.Item(18,4).Value = "COMERCIO CENTRAL "
$ws.Cells.Item(18,5).Value = 1785

$ws.Cells.Item(19,1).Value = 44512
$ws.Cells.Item(19,4).Value = "COMERCIO CENTRAL "
$ws.Cells.Item(19,5).Value = 13805

$ws.Cells.Item(20,1).Value = 44513
$ws.Cells.Item(20,4).Value = "CANCELADA"
$ws.Cells.Item(20,4).Font.Bold = $true
$ws.Cells.Item(20,4).Font.Size = 12
$ws.Cells.Item(20,4).Font.Color = 255
$ws.Cells.Item(20,5).Value = 0

$ws.Cells.Item(21,1).Value = 44513
$ws.Cells.Item(21,4).Value = "OBRADOR"
$ws.Cells.Item(21,5).Value = 18875

$ws.Cells.Item(22,1).Value = 44513
$ws.Cells.Item(22,4).Value = "OBRADOR"
$ws.Cells.Item(22,5).Value = 10476

$ws.Cells.Item(23,1).Value = 44513
$ws.Cells.Item(23,4).Value = "CANCELADA"
$ws.Cells.Item(23,5).Value = 0

$ws.Cells.Item(24,1).Value = 44513
$ws.Cells.Item(24,4).Value = "COMERCIO CENTRAL "
$ws.Cells.Item(24,5).Value = 219644

$ws.Cells.Item(25,1).Value = 44513
$ws.Cells.Item(25,4).Value = "COMERCIO CENTRAL "
$ws.Cells.Item(25,5).Value = 2546

$ws.Cells.Item(26,1).Value = 44513
$ws.Cells.Item(26,4).Value = "ABASTOS DE 4 CARNES 11 SUR "
$ws.Cells.Item(26,5).Value = 39216

# Move the view so that row 10 is at the top and D27 is the active selection,
# matching where the user continued entering data.
$ws.Range("D27").Select()
$excel.ActiveWindow.ScrollRow = 10
$excel.ActiveWindow.ScrollColumn = 1
